# Weekly update for "Hortaliza, Vega Monumental Concepción - Zanahoria":
# a new week's pair of records (Primera / Segunda) is inserted at the top
# of the dated block (rows 131-132), pushing the existing rows 131-166
# down to 133-168. Fill the two freshly inserted rows with the new
# week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("131:132").Insert()

# Row 131 - Zanahoria, Primera, semana del 2021-12-29
$ws.Cells.Item(131, 1).Value = 11
$ws.Cells.Item(131, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(131, 3).Value = "Bíobío"
$ws.Cells.Item(131, 4).Value = 44559
$ws.Cells.Item(131, 5).Value = 8
$ws.Cells.Item(131, 6).Value = 100114013
$ws.Cells.Item(131, 7).Value = "Zanahoria"
$ws.Cells.Item(131, 8).Value = "Sin especificar"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 500
$ws.Cells.Item(131, 11).Value = 7000
$ws.Cells.Item(131, 12).Value = 8000
$ws.Cells.Item(131, 13).Value = 7400
$ws.Cells.Item(131, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(131, 15).Value = "Región de Ñuble"
$ws.Cells.Item(131, 16).Value = 370
$ws.Cells.Item(131, 17).Value = 20
$ws.Cells.Item(131, 18).Value = "Hortaliza"

# Row 132 - Zanahoria, Segunda, semana del 2021-12-29
$ws.Cells.Item(132, 1).Value = 11
$ws.Cells.Item(132, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(132, 3).Value = "Bíobío"
$ws.Cells.Item(132, 4).Value = 44559
$ws.Cells.Item(132, 5).Value = 8
$ws.Cells.Item(132, 6).Value = 100114013
$ws.Cells.Item(132, 7).Value = "Zanahoria"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Segunda"
$ws.Cells.Item(132, 10).Value = 300
$ws.Cells.Item(132, 11).Value = 6000
$ws.Cells.Item(132, 12).Value = 6000
$ws.Cells.Item(132, 13).Value = 6000
$ws.Cells.Item(132, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(132, 15).Value = "Región de Ñuble"
$ws.Cells.Item(132, 16).Value = 300
$ws.Cells.Item(132, 17).Value = 20
$ws.Cells.Item(132, 18).Value = "Hortaliza"
